$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Single-value cells that just get their number text swapped out.
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "208"
$t.Cell(6, 1).Range.Text  = "0.00016"
$t.Cell(7, 1).Range.Text  = "0.00007"
$t.Cell(9, 1).Range.Text  = "0.00016"
$t.Cell(10, 1).Range.Text = "0.00016"
$t.Cell(11, 1).Range.Text = "0.00016"
$t.Cell(12, 1).Range.Text = "0.00767"

# Rows 44-46 previously held a whole tab-separated line of stats crammed
# into one run; they collapse down to just the leading count column.
$t.Cell(44, 1).Range.Text = "100"
$t.Cell(45, 1).Range.Text = "0.01"
$t.Cell(46, 1).Range.Text = "362"
